$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Result" slide): title + content tweaks
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)

# Title: "Result" -> "Result/Future Outlook"
$titleShape = $slide10.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Result/Future Outlook"

# Content: append a new sentence to the existing paragraph text.
$contentShape = $slide10.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange
$contentRange.Text = $contentRange.Text + " I would recommend possibly making a dataset on players playing for their country and making an average based on their performance to see which teams has the most elite players."

# ---------------------------------------------------------------------------
# Slide 3 ("Problems" slide): split the trailing run into three runs and
# rewrite the middle part ("scored, I used left outer join to combine." ->
# "scored.The").
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$problemShape = $slide3.Shapes.Item(2)
$problemRange = $problemShape.TextFrame.TextRange

$oldTail = " and number of goals scored, I used left outer join to combine. The other excel sheet had a large array of data"
$fullText = $problemRange.Text
$tailStart0 = $fullText.IndexOf($oldTail)

if ($tailStart0 -ge 0) {
    $part1 = " and number of goals "
    $part2 = "scored.The"
    $part3 = " other excel sheet had a large array of data"

    $tailStart1 = $tailStart0 + 1
    $problemRange.Characters($tailStart1, $oldTail.Length).Text = ($part1 + $part2 + $part3)

    $p1Start = $tailStart1
    $p2Start = $p1Start + $part1.Length
    $p3Start = $p2Start + $part2.Length

    $problemRange.Characters($p1Start, $part1.Length).Text = $part1
    $problemRange.Characters($p2Start, $part2.Length).Text = $part2
    $problemRange.Characters($p3Start, $part3.Length).Text = $part3
}

# ---------------------------------------------------------------------------
# Slide 4 (Excel snippet picture slide): give the empty title placeholder
# centered text "Snippet of Excel ".
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$snippetTitle = $slide4.Shapes.Item(1)
$snippetRange = $snippetTitle.TextFrame.TextRange
$snippetRange.Text = "Snippet of Excel "
$snippetRange.ParagraphFormat.Alignment = 2
